$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.807.98'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.789.42'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.13%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.09'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.11%  '

$ws.Range("E6").Value = '  -0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5124'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3884'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07814'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -7.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.091'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.77'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.80%  '

$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.213'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.18'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.219'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.777.96'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.45'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.63%  '

$ws.Range("E18").Value = '  -4.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06519'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.03'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.910'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.889.38'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.00'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.226'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.33'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.26'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.988.63'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.355'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.83%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.87'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1075'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.037'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.621'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.485'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07030'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02303'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.730'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2126'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.51'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.991'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6083'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.148'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.89%  '

$ws.Range("E44").Value = '  -6.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.10'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5890'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.695'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.34'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.199'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.909'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06809'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.53%  '
